# MOSFERATU JLCPCB BOM.xlsx - designator reassignment
# R23 moves from the "1kΩ" row to the "150kΩ" row, R24 moves from the
# "1kΩ" row to the "10kΩ" row, and R4 is added to the "4.7kΩ" row
# (alongside RLEDFX).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = "'R2A, R3A, R10, R14, R20, R24"
$ws.Range("B21").Value = "'R15, R23"
$ws.Range("B23").Value = "'R3, R13, R17, R22, R36"
$ws.Range("B30").Value = "'R4, RLEDFX"

# Match the author's final on-screen selection/scroll state.
[void]$ws.Range("B28").Select()
